$d = $word.ActiveDocument

# 1. "florissent des premiers &" -> "florissent des premiers, &"
$d.Content.Find.Execute("florissent des premiers &", $true, $false, $false, $false, $false,
                         $true, 1, $false, "florissent des premiers, &", 2)

# 2. "florissant plus tard les fleurs viennent" -> "florissant plus tard, les fleurs viennent"
$d.Content.Find.Execute("florissant plus tard les fleurs viennent", $true, $false, $false, $false, $false,
                         $true, 1, $false, "florissant plus tard, les fleurs viennent", 2)

# 3. "</pl>" -> "</pl>," (literal angle-bracket tag text inside the run)
$d.Content.Find.Execute("</pl>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "</pl>,", 2)

# 4. "hiver &" -> "hiver, &" (only the "l'hiver & l'esté" occurrence matches this text;
#    the other "hiver" in the document reads "hiver affin que" and is left untouched)
$d.Content.Find.Execute("hiver &", $true, $false, $false, $false, $false,
                         $true, 1, $false, "hiver, &", 2)
